# Refresh the cryptocurrency price/volume figures on Sheet1 (columns D and E)
# to match the latest scrape, per the GitHub Actions update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '26.658.17'
$ws.Range('E2').Value = '  +0.57%  '
$ws.Range('D3').Value = '1.598.28'
$ws.Range('E3').Value = '  +1.04%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = "'210.88"
$ws.Range('E5').Value = '  +0.07%  '
$ws.Range('D6').Value = "'0.513"
$ws.Range('E6').Value = '  +1.40%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').Value = "'0.0617"
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -1.16%  '
$ws.Range('D10').Value = "'19.61"
$ws.Range('E10').Value = '  +0.71%  '
$ws.Range('D11').Value = "'0.0842"
$ws.Range('E11').Value = '  +1.26%  '
$ws.Range('D12').Value = '1.821.99'
$ws.Range('E12').Value = '  +1.05%  '
$ws.Range('D13').Value = '1.583.43'
$ws.Range('E13').Value = '  +0.10%  '
$ws.Range('D14').Value = "'4.02"
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('D15').Value = "'0.521"
$ws.Range('E15').Value = '  -1.34%  '
$ws.Range('D16').Value = "'64.75"
$ws.Range('E16').Value = '  +1.21%  '
$ws.Range('D17').Value = '26.627.56'
$ws.Range('E17').Value = '  +0.36%  '
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('D19').Value = "'208.72"
$ws.Range('E19').Value = '  +0.29%  '
$ws.Range('E20').Value = '  -0.11%  '
$ws.Range('D21').Value = "'6.75"
$ws.Range('E21').Value = '  +1.09%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').Value = "'2.30"
$ws.Range('E23').Value = '  -3.59%  '
$ws.Range('D24').Value = "'8.89"
$ws.Range('E24').Value = '  +0.11%  '
$ws.Range('D25').Value = "'145.59"
$ws.Range('E25').Value = '  -0.60%  '
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').Value = "'7.24"
$ws.Range('E27').Value = '  -2.42%  '
$ws.Range('E28').Value = '  +2.25%  '
$ws.Range('D29').Value = "'15.24"
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('E30').Value = '  +0.82%  '
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('E32').Value = '  -0.68%  '
$ws.Range('D33').Value = "'0.652"
$ws.Range('E33').Value = '  -0.43%  '
$ws.Range('E34').Value = '  -0.30%  '
$ws.Range('D35').Value = '1.291.44'
$ws.Range('E35').Value = '  -1.04%  '
$ws.Range('E36').Value = '  +0.47%  '
$ws.Range('E37').Value = '  -1.26%  '
$ws.Range('E38').Value = '  -0.84%  '
$ws.Range('D39').Value = "'0.845"
$ws.Range('E39').Value = '  +3.13%  '
$ws.Range('E40').Value = '  -0.11%  '
$ws.Range('E41').Value = '  +2.00%  '
$ws.Range('E42').Value = '  +1.62%  '
$ws.Range('E43').Value = '  +0.62%  '
$ws.Range('D44').Value = "'63.82"
$ws.Range('D45').Value = '1.734.55'
$ws.Range('E45').Value = '  +1.04%  '
$ws.Range('D46').Value = "'0.892"
$ws.Range('E46').Value = '  +7.52%  '
$ws.Range('D47').Value = "'90.08"
$ws.Range('E47').Value = '  +1.36%  '
$ws.Range('E48').Value = '  -0.38%  '
$ws.Range('D49').Value = "'0.100"
$ws.Range('E49').Value = '  +2.13%  '
$ws.Range('D50').Value = "'0.0505"
$ws.Range('E50').Value = '  -0.24%  '
$ws.Range('D51').Value = "'7.48"
$ws.Range('E51').Value = '  +0.47%  '
